# Auto update Excel log
# Appends newly-logged sensor events (2026-01-30, ~16:10-16:11) to the
# PIR, Humidity, Proximity and Camera sheets, matching the source data
# exporter's behaviour. Values are entered with a leading apostrophe so
# Excel stores them as literal text (matching the existing inlineStr
# cells) instead of auto-converting dates/times/percentages to numbers.

$wb = $excel.ActiveWorkbook

# ---- PIR sheet: append rows 84..97 ----
$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A84").Value = "'2026-01-30"
$ws.Range("B84").Value = "'16:10:23"
$ws.Range("C84").Value = "'16:00"
$ws.Range("D84").Value = "'Bathroom"
$ws.Range("E84").Value = "'No Motion"
$ws.Range("F84").Value = "'Inactive"
$ws.Range("A85").Value = "'2026-01-30"
$ws.Range("B85").Value = "'16:10:23"
$ws.Range("C85").Value = "'16:00"
$ws.Range("D85").Value = "'Bathroom"
$ws.Range("E85").Value = "'No Motion"
$ws.Range("F85").Value = "'Inactive"
$ws.Range("A86").Value = "'2026-01-30"
$ws.Range("B86").Value = "'16:10:27"
$ws.Range("C86").Value = "'16:00"
$ws.Range("D86").Value = "'Bathroom"
$ws.Range("E86").Value = "'No Motion"
$ws.Range("F86").Value = "'Inactive"
$ws.Range("A87").Value = "'2026-01-30"
$ws.Range("B87").Value = "'16:10:33"
$ws.Range("C87").Value = "'16:00"
$ws.Range("D87").Value = "'Bathroom"
$ws.Range("E87").Value = "'No Motion"
$ws.Range("F87").Value = "'Inactive"
$ws.Range("A88").Value = "'2026-01-30"
$ws.Range("B88").Value = "'16:10:37"
$ws.Range("C88").Value = "'16:00"
$ws.Range("D88").Value = "'Bathroom"
$ws.Range("E88").Value = "'No Motion"
$ws.Range("F88").Value = "'Inactive"
$ws.Range("A89").Value = "'2026-01-30"
$ws.Range("B89").Value = "'16:10:42"
$ws.Range("C89").Value = "'16:00"
$ws.Range("D89").Value = "'Bathroom"
$ws.Range("E89").Value = "'No Motion"
$ws.Range("F89").Value = "'Inactive"
$ws.Range("A90").Value = "'2026-01-30"
$ws.Range("B90").Value = "'16:10:47"
$ws.Range("C90").Value = "'16:00"
$ws.Range("D90").Value = "'Bathroom"
$ws.Range("E90").Value = "'No Motion"
$ws.Range("F90").Value = "'Inactive"
$ws.Range("A91").Value = "'2026-01-30"
$ws.Range("B91").Value = "'16:10:53"
$ws.Range("C91").Value = "'16:00"
$ws.Range("D91").Value = "'Bathroom"
$ws.Range("E91").Value = "'No Motion"
$ws.Range("F91").Value = "'Inactive"
$ws.Range("A92").Value = "'2026-01-30"
$ws.Range("B92").Value = "'16:10:57"
$ws.Range("C92").Value = "'16:00"
$ws.Range("D92").Value = "'Bathroom"
$ws.Range("E92").Value = "'No Motion"
$ws.Range("F92").Value = "'Inactive"
$ws.Range("A93").Value = "'2026-01-30"
$ws.Range("B93").Value = "'16:11:03"
$ws.Range("C93").Value = "'16:00"
$ws.Range("D93").Value = "'Bathroom"
$ws.Range("E93").Value = "'No Motion"
$ws.Range("F93").Value = "'Inactive"
$ws.Range("A94").Value = "'2026-01-30"
$ws.Range("B94").Value = "'16:11:07"
$ws.Range("C94").Value = "'16:00"
$ws.Range("D94").Value = "'Bathroom"
$ws.Range("E94").Value = "'No Motion"
$ws.Range("F94").Value = "'Inactive"
$ws.Range("A95").Value = "'2026-01-30"
$ws.Range("B95").Value = "'16:11:12"
$ws.Range("C95").Value = "'16:00"
$ws.Range("D95").Value = "'Bathroom"
$ws.Range("E95").Value = "'No Motion"
$ws.Range("F95").Value = "'Inactive"
$ws.Range("A96").Value = "'2026-01-30"
$ws.Range("B96").Value = "'16:11:17"
$ws.Range("C96").Value = "'16:00"
$ws.Range("D96").Value = "'Bathroom"
$ws.Range("E96").Value = "'No Motion"
$ws.Range("F96").Value = "'Inactive"
$ws.Range("A97").Value = "'2026-01-30"
$ws.Range("B97").Value = "'16:11:23"
$ws.Range("C97").Value = "'16:00"
$ws.Range("D97").Value = "'Bathroom"
$ws.Range("E97").Value = "'No Motion"
$ws.Range("F97").Value = "'Inactive"

# ---- Humidity sheet: append rows 68..76 ----
$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A68").Value = "'2026-01-30"
$ws.Range("B68").Value = "'16:10:23"
$ws.Range("C68").Value = "'16:00"
$ws.Range("D68").Value = "'Bathroom"
$ws.Range("E68").Value = "'87.5%"
$ws.Range("F68").Value = "'Active"
$ws.Range("A69").Value = "'2026-01-30"
$ws.Range("B69").Value = "'16:10:23"
$ws.Range("C69").Value = "'16:00"
$ws.Range("D69").Value = "'Bathroom"
$ws.Range("E69").Value = "'86.5%"
$ws.Range("F69").Value = "'Active"
$ws.Range("A70").Value = "'2026-01-30"
$ws.Range("B70").Value = "'16:10:28"
$ws.Range("C70").Value = "'16:00"
$ws.Range("D70").Value = "'Bathroom"
$ws.Range("E70").Value = "'87.4%"
$ws.Range("F70").Value = "'Active"
$ws.Range("A71").Value = "'2026-01-30"
$ws.Range("B71").Value = "'16:10:33"
$ws.Range("C71").Value = "'16:00"
$ws.Range("D71").Value = "'Bathroom"
$ws.Range("E71").Value = "'86.6%"
$ws.Range("F71").Value = "'Active"
$ws.Range("A72").Value = "'2026-01-30"
$ws.Range("B72").Value = "'16:10:38"
$ws.Range("C72").Value = "'16:00"
$ws.Range("D72").Value = "'Bathroom"
$ws.Range("E72").Value = "'87.5%"
$ws.Range("F72").Value = "'Active"
$ws.Range("A73").Value = "'2026-01-30"
$ws.Range("B73").Value = "'16:10:48"
$ws.Range("C73").Value = "'16:00"
$ws.Range("D73").Value = "'Bathroom"
$ws.Range("E73").Value = "'87.6%"
$ws.Range("F73").Value = "'Active"
$ws.Range("A74").Value = "'2026-01-30"
$ws.Range("B74").Value = "'16:10:58"
$ws.Range("C74").Value = "'16:00"
$ws.Range("D74").Value = "'Bathroom"
$ws.Range("E74").Value = "'87.5%"
$ws.Range("F74").Value = "'Active"
$ws.Range("A75").Value = "'2026-01-30"
$ws.Range("B75").Value = "'16:11:08"
$ws.Range("C75").Value = "'16:00"
$ws.Range("D75").Value = "'Bathroom"
$ws.Range("E75").Value = "'87.5%"
$ws.Range("F75").Value = "'Active"
$ws.Range("A76").Value = "'2026-01-30"
$ws.Range("B76").Value = "'16:11:18"
$ws.Range("C76").Value = "'16:00"
$ws.Range("D76").Value = "'Bathroom"
$ws.Range("E76").Value = "'87.5%"
$ws.Range("F76").Value = "'Active"

# ---- Proximity sheet: append rows 17..19 ----
$ws = $wb.Worksheets.Item("Proximity")
$ws.Range("A17").Value = "'2026-01-30"
$ws.Range("B17").Value = "'16:10:32"
$ws.Range("C17").Value = "'16:00"
$ws.Range("D17").Value = "'Living Room Main Door"
$ws.Range("E17").Value = "'EXIT"
$ws.Range("F17").Value = "'User EXITED Living Room Main Door"
$ws.Range("A18").Value = "'2026-01-30"
$ws.Range("B18").Value = "'16:10:35"
$ws.Range("C18").Value = "'16:00"
$ws.Range("D18").Value = "'Living Room Main Door"
$ws.Range("E18").Value = "'ENTER"
$ws.Range("F18").Value = "'User ENTERED Living Room Main Door"
$ws.Range("A19").Value = "'2026-01-30"
$ws.Range("B19").Value = "'16:11:17"
$ws.Range("C19").Value = "'16:00"
$ws.Range("D19").Value = "'Living Room Main Door"
$ws.Range("E19").Value = "'EXIT"
$ws.Range("F19").Value = "'User EXITED Living Room Main Door"

# ---- Camera sheet: append rows 17..19 ----
$ws = $wb.Worksheets.Item("Camera")
$ws.Range("A17").Value = "'2026-01-30"
$ws.Range("B17").Value = "'16:10:32"
$ws.Range("C17").Value = "'16:00"
$ws.Range("D17").Value = "'Living Room Main Door"
$ws.Range("E17").Value = "'Image Captured (EXIT)"
$ws.Range("F17").Value = "'Active"
$ws.Range("A18").Value = "'2026-01-30"
$ws.Range("B18").Value = "'16:10:35"
$ws.Range("C18").Value = "'16:00"
$ws.Range("D18").Value = "'Living Room Main Door"
$ws.Range("E18").Value = "'Image Captured (ENTER)"
$ws.Range("F18").Value = "'Active"
$ws.Range("A19").Value = "'2026-01-30"
$ws.Range("B19").Value = "'16:11:17"
$ws.Range("C19").Value = "'16:00"
$ws.Range("D19").Value = "'Living Room Main Door"
$ws.Range("E19").Value = "'Image Captured (EXIT)"
$ws.Range("F19").Value = "'Active"
Write-Host "Appended new sensor log rows to PIR, Humidity, Proximity and Camera sheets."
